# Update crypto price/volume columns (D, E) for rows 2-51 as plain text values,
# matching the commit's GitHub Actions refresh. Each value is written via a scratch
# formula cell + copy/paste-values so Excel keeps it as text (no numeric coercion,
# no NumberFormat/style changes), mirroring the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

$scratch.Formula = "=""38.208.94"""
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.57%  """
$scratch.Copy()
$ws.Range("E2").PasteSpecial(-4163)

$scratch.Formula = "=""2.061.09"""
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.50%  """
$scratch.Copy()
$ws.Range("E3").PasteSpecial(-4163)

$scratch.Formula = "=""  +0.06%  """
$scratch.Copy()
$ws.Range("E4").PasteSpecial(-4163)

$scratch.Formula = "=""230.66"""
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.04%  """
$scratch.Copy()
$ws.Range("E5").PasteSpecial(-4163)

$scratch.Formula = "=""  +2.20%  """
$scratch.Copy()
$ws.Range("E6").PasteSpecial(-4163)

$scratch.Formula = "=""58.13"""
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)

$scratch.Formula = "=""  +7.78%  """
$scratch.Copy()
$ws.Range("E7").PasteSpecial(-4163)

$scratch.Formula = "=""1.00"""
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)

$scratch.Formula = "=""  -0.06%  """
$scratch.Copy()
$ws.Range("E8").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.63%  """
$scratch.Copy()
$ws.Range("E9").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.91%  """
$scratch.Copy()
$ws.Range("E10").PasteSpecial(-4163)

$scratch.Formula = "=""  -0.49%  """
$scratch.Copy()
$ws.Range("E11").PasteSpecial(-4163)

$scratch.Formula = "=""2.369.29"""
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.31%  """
$scratch.Copy()
$ws.Range("E12").PasteSpecial(-4163)

$scratch.Formula = "=""  +4.43%  """
$scratch.Copy()
$ws.Range("E13").PasteSpecial(-4163)

$scratch.Formula = "=""20.72"""
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.89%  """
$scratch.Copy()
$ws.Range("E14").PasteSpecial(-4163)

$scratch.Formula = "=""0.757"""
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.28%  """
$scratch.Copy()
$ws.Range("E15").PasteSpecial(-4163)

$scratch.Formula = "=""5.30"""
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)

$scratch.Formula = "=""  +4.61%  """
$scratch.Copy()
$ws.Range("E16").PasteSpecial(-4163)

$scratch.Formula = "=""2.062.07"""
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.59%  """
$scratch.Copy()
$ws.Range("E17").PasteSpecial(-4163)

$scratch.Formula = "=""38.120.86"""
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.65%  """
$scratch.Copy()
$ws.Range("E18").PasteSpecial(-4163)

$scratch.Formula = "=""6.18"""
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)

$scratch.Formula = "=""  +2.47%  """
$scratch.Copy()
$ws.Range("E19").PasteSpecial(-4163)

$scratch.Formula = "=""69.87"""
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)

$scratch.Formula = "=""  +2.34%  """
$scratch.Copy()
$ws.Range("E20").PasteSpecial(-4163)

$scratch.Formula = "=""0.0₃0833"""
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.17%  """
$scratch.Copy()
$ws.Range("E21").PasteSpecial(-4163)

$scratch.Formula = "=""225.10"""
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)

$scratch.Formula = "=""  +1.78%  """
$scratch.Copy()
$ws.Range("E22").PasteSpecial(-4163)

$scratch.Formula = "=""  +0.09%  """
$scratch.Copy()
$ws.Range("E23").PasteSpecial(-4163)

$scratch.Formula = "=""  +1.63%  """
$scratch.Copy()
$ws.Range("E24").PasteSpecial(-4163)

$scratch.Formula = "=""  +4.40%  """
$scratch.Copy()
$ws.Range("E25").PasteSpecial(-4163)

$scratch.Formula = "=""9.32"""
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.11%  """
$scratch.Copy()
$ws.Range("E26").PasteSpecial(-4163)

$scratch.Formula = "=""166.34"""
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)

$scratch.Formula = "=""  +0.58%  """
$scratch.Copy()
$ws.Range("E27").PasteSpecial(-4163)

$scratch.Formula = "=""0.135"""
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)

$scratch.Formula = "=""  +8.94%  """
$scratch.Copy()
$ws.Range("E28").PasteSpecial(-4163)

$scratch.Formula = "=""19.07"""
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.24%  """
$scratch.Copy()
$ws.Range("E29").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.17%  """
$scratch.Copy()
$ws.Range("E30").PasteSpecial(-4163)

$scratch.Formula = "=""  +2.23%  """
$scratch.Copy()
$ws.Range("E31").PasteSpecial(-4163)

$scratch.Formula = "=""  +2.84%  """
$scratch.Copy()
$ws.Range("E32").PasteSpecial(-4163)

$scratch.Formula = "=""  +5.83%  """
$scratch.Copy()
$ws.Range("E33").PasteSpecial(-4163)

$scratch.Formula = "=""0.0616"""
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)

$scratch.Formula = "=""  +1.92%  """
$scratch.Copy()
$ws.Range("E34").PasteSpecial(-4163)

$scratch.Formula = "=""  +7.21%  """
$scratch.Copy()
$ws.Range("E35").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.53%  """
$scratch.Copy()
$ws.Range("E36").PasteSpecial(-4163)

$scratch.Formula = "=""6.10"""
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)

$scratch.Formula = "=""  +15.85%  """
$scratch.Copy()
$ws.Range("E37").PasteSpecial(-4163)

$scratch.Formula = "=""3.34"""
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)

$scratch.Formula = "=""  +7.48%  """
$scratch.Copy()
$ws.Range("E38").PasteSpecial(-4163)

$scratch.Formula = "=""  -0.09%  """
$scratch.Copy()
$ws.Range("E39").PasteSpecial(-4163)

$scratch.Formula = "=""98.52"""
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)

$scratch.Formula = "=""  +5.15%  """
$scratch.Copy()
$ws.Range("E40").PasteSpecial(-4163)

$scratch.Formula = "=""  +2.42%  """
$scratch.Copy()
$ws.Range("E41").PasteSpecial(-4163)

$scratch.Formula = "=""1.486.38"""
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)

$scratch.Formula = "=""  +1.46%  """
$scratch.Copy()
$ws.Range("E42").PasteSpecial(-4163)

$scratch.Formula = "=""16.94"""
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)

$scratch.Formula = "=""  +4.27%  """
$scratch.Copy()
$ws.Range("E43").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.33%  """
$scratch.Copy()
$ws.Range("E44").PasteSpecial(-4163)

$scratch.Formula = "=""  +3.96%  """
$scratch.Copy()
$ws.Range("E45").PasteSpecial(-4163)

$scratch.Formula = "=""  +1.55%  """
$scratch.Copy()
$ws.Range("E46").PasteSpecial(-4163)

$scratch.Formula = "=""4.10"""
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)

$scratch.Formula = "=""  +19.15%  """
$scratch.Copy()
$ws.Range("E47").PasteSpecial(-4163)

$scratch.Formula = "=""  +2.63%  """
$scratch.Copy()
$ws.Range("E48").PasteSpecial(-4163)

$scratch.Formula = "=""2.96"""
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)

$scratch.Formula = "=""  +2.68%  """
$scratch.Copy()
$ws.Range("E49").PasteSpecial(-4163)

$scratch.Formula = "=""  +0.15%  """
$scratch.Copy()
$ws.Range("E50").PasteSpecial(-4163)

$scratch.Formula = "=""2.252.52"""
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)

$scratch.Formula = "=""  +2.63%  """
$scratch.Copy()
$ws.Range("E51").PasteSpecial(-4163)

$scratch.Clear()
$excel.CutCopyMode = 0
